# Fix 'or', 'redirect' words
#
# 1) "...actor's request and retrieve the data from User." -> "...request or retrieve..."
# 2) Every "RedirectRequest" -> "RefreshRequest" (5 occurrences across the
#    two sequence-diagram tables plus the final summary table).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "and" -> "or" in the "Prepare a query ..." cell.
#    Split the run on the word boundaries (mirrors how Word splits a run
#    when you retype a selected word) so the surrounding text keeps its
#    original formatting and only the new "or" run is distinct.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("request and retrieve the data from User.")
if ($found) {
    $start = $rng.Start
    # "request " is 8 characters, "and" is the next 3.
    $andRange = $d.Range($start + 8, $start + 11)
    $andRange.Text = "or"

    # Force a clean run boundary around the freshly-typed "or" so it does
    # not simply get swallowed back into the neighbouring runs.
    $orRange = $d.Range($start + 8, $start + 10)
    $orRange.Font.Bold = $true
    $orRange.Font.Bold = $false
} else {
    Write-Host "WARNING: 'request and retrieve the data from User.' not found"
}

# ---------------------------------------------------------------------
# 2) RedirectRequest -> RefreshRequest (5 occurrences).
#    Each one is spelled "R" + "edirectRequest" (two runs) except the
#    "UserOperator generates RedirectRequest" cell, where the whole
#    phrase is a single run. A plain targeted replace of the
#    "edirectRequest" substring keeps every other run (and the leading
#    "R" run, when present) untouched.
# ---------------------------------------------------------------------
$count = 0
while ($true) {
    $rng = $d.Content
    $found = $rng.Find.Execute("edirectRequest")
    if (-not $found) { break }
    $rng.Text = "efreshRequest"
    $count = $count + 1
    if ($count -gt 20) { break }
}
Write-Host "Replaced $count occurrence(s) of 'edirectRequest' with 'efreshRequest'"
